$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the numeric "شماره ملی" (national ID) values in G2:G7 with their
# correct text-formatted national-ID numbers (stored as shared strings,
# matching the column's existing Text number format).
$ws.Range("G2").Value = "3100471249"
$ws.Range("G3").Value = "3378946117"
$ws.Range("G4").Value = "8796532154"
$ws.Range("G5").Value = "3657420985"
$ws.Range("G6").Value = "3935895853"
$ws.Range("G7").Value = "8796532154"

# Update the sheet's current selection to match the edited range.
$ws.Range("G2:G7").Select()
